# "Generate Report for Handoff"
#
# This updates the localization-status report:
#   - Status cells move from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The "Latest HO Xliff Generate Date" / "Latest Handback DateTime" timestamps
#     (that shared the same value as the status datetime) are refreshed
#   - The zh-cn sheet's "Latest Handoff Datetime" is refreshed
#   - The now-shorter status text no longer needs as wide a column, so the
#     columns that show the status get narrowed on all three sheets

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview"
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$dede     = $wb.Worksheets.Item(3)   # "de-de"

# --- Overview sheet ---
# E2 (zh-cn status) and F2 (de-de status)
$overview.Range("E2:F2").Value = "Ready for handoff"
# G2 "Latest HO Xliff Generate Date"
$overview.Range("G2").Value = "2016-09-04 01:05:18"

# --- zh-cn sheet ---
# C2 Status
$zhcn.Range("C2").Value = "Ready for handoff"
# H2 "Latest Handoff Datetime"
$zhcn.Range("H2").Value = "2016-09-04 01:05:13"

# --- de-de sheet ---
# C2 Status
$dede.Range("C2").Value = "Ready for handoff"
# H2 "Latest Handoff Datetime" (shares value with Overview's G2 date)
$dede.Range("H2").Value = "2016-09-04 01:05:18"

# --- Narrow the status columns now that the text is shorter ---
# Overview columns E and F
$overview.Range("E1:F1").ColumnWidth = 16.33
# zh-cn / de-de column C
$zhcn.Range("C1").ColumnWidth = 16.33
$dede.Range("C1").ColumnWidth = 16.33
